$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.462.17"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.800.30"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.39"
$ws.Range("E5").Value = "  +4.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.57"
$ws.Range("E6").Value = "  -3.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("E7").Value = "  +4.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  +2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.18"
$ws.Range("E10").Value = "  -2.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("E11").Value = "  -1.92%  "

$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.83"
$ws.Range("E13").Value = "  -2.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.68"
$ws.Range("E14").Value = "  +0.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.244.17"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.790.12"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.883"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.245.55"
$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("E19").Value = "  +7.04%  "

$ws.Range("E20").Value = "  -5.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").Value = "  -2.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0989"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.66"
$ws.Range("E23").Value = "  -3.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.37"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("E30").Value = "  -2.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "50.48"
$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.81"
$ws.Range("E32").Value = "  -3.91%  "

$ws.Range("E33").Value = "  +3.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0440"
$ws.Range("E34").Value = "  +24.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0819"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.97"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.05"
$ws.Range("E38").Value = "  -2.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.18"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.97"
$ws.Range("E40").Value = "  -6.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.60"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.84"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("E44").Value = "  +1.34%  "

$ws.Range("E45").Value = "  -1.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.069.46"
$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  +2.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.61"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("E50").Value = "  +4.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.44"
$ws.Range("E51").Value = "  -0.16%  "
